$wb = $excel.ActiveWorkbook

# Style used by the existing "Merkmalswerte" header row (bold, bordered,
# centered/top-aligned) -- reuse it for the new sheets' header rows so we
# don't introduce any new style records.
$headerStyleRow = $wb.Worksheets.Item(1).Range("A1:G1")

# --- New sheet "Merkmale" (added right after "Merkmalswerte") ---
$wsMerkmale = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsMerkmale.Name = "Merkmale"
$wsMerkmale.Range("A1").Value = "ATINN"
$wsMerkmale.Range("B1").Value = "ADZHL"
$wsMerkmale.Range("C1").Value = "Merkmalsname"
$wsMerkmale.Range("D1").Value = "Merkmalsbezeichnung (deutsch)"
$wsMerkmale.Range("E1").Value = "Langtext"
$wsMerkmale.Range("F1").Value = "Merkmalsbezeichnung (englisch)"
$wsMerkmale.Range("G1").Value = "Langtext"

$headerStyleRow.Copy()
$wsMerkmale.Range("A1:G1").PasteSpecial(-4122)

# --- New sheet "Konditionen" (added after "Merkmale") ---
$wsKonditionen = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsKonditionen.Name = "Konditionen"
$wsKonditionen.Range("A2").Value = "Kondition"
$wsKonditionen.Range("B2").Value = "deutsch"
$wsKonditionen.Range("C2").Value = "englisch"

$headerStyleRow.Range("A1:C1").Copy()
$wsKonditionen.Range("A2:C2").PasteSpecial(-4122)
